# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text (row 1)
$ws.Range("A1").Value = "Datos actualizados a 29 de Marzo de 2020 a las 12:20"

# Madrid (row 4) - updated totals
$ws.Range("B4").Value = 22677
$ws.Range("C4").Value = 7491
$ws.Range("D4").Value = 12104
$ws.Range("E4").Value = 3082

# Navarra (row 9) - updated totals
$ws.Range("B9").Value = 2011
$ws.Range("C9").Value = 125
$ws.Range("D9").Value = 1802
$ws.Range("E9").Value = 84

# Cantabria's case count grew enough to move above Pontevedra and Caceres
# in this ranking table, so rows 21-23 shift accordingly.
# Row 21: now Cantabria (was Pontevedra)
$ws.Range("A21").Value = "Cantabria"
$ws.Range("B21").Value = 1023
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 972
$ws.Range("E21").Value = 26

# Row 22: now Pontevedra (was Caceres)
$ws.Range("A22").Value = "Pontevedra"
$ws.Range("B22").Value = 960
$ws.Range("C22").Value = 95
$ws.Range("D22").Value = 923
$ws.Range("E22").Value = 9

# Row 23: now Caceres (was Cantabria)
$ws.Range("A23").Value = "Caceres"
$ws.Range("B23").Value = 957
$ws.Range("C23").Value = 38
$ws.Range("D23").Value = 841
$ws.Range("E23").Value = 78

# Melilla (row 54) - updated totals
$ws.Range("B54").Value = 48
$ws.Range("D54").Value = 47
